$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'curvy leggings'
    2 = 'custom compression pants'
    3 = 'custom ultimate frisbee'
    4 = 'cut leggings for women'
    5 = 'cw amazon channels'
    6 = 'cw arrow costume'
    7 = 'cw designs'
    8 = 'cw everyday plus'
    9 = 'cw flash costume'
    10 = 'cw flash costume adult'
    11 = 'cw flash season 4'
    12 = 'cw key'
    13 = 'cw live'
    14 = 'cw merchandise'
    15 = 'cw now'
    16 = 'cw post apparel'
    17 = 'cw shirt'
    18 = 'cw simply extra'
    19 = 'cw straight key'
    20 = 'cw trainer'
    21 = 'cw x'
    22 = 'cw x bra'
    23 = 'cw x compression'
    24 = 'cycle apparel women'
    25 = 'cycle leggings women'
    26 = 'cycle pants women'
    27 = 'cycle tights women'
    28 = 'cycle wear women'
    29 = 'cycling 3 4 pants women'
    30 = 'cycling art'
    31 = 'cycling capri'
    32 = 'cycling capri pants'
    33 = 'cycling capris'
    34 = 'cycling capris for women'
    35 = 'cycling compression'
    36 = 'cycling compression shorts'
    37 = 'cycling compression tights'
    38 = 'cycling knee brace'
    39 = 'cycling leg warmers'
    40 = 'cycling leg warmers womens'
    41 = 'cycling leggings women'
    42 = 'cycling leggins for men'
    43 = 'cycling long pants'
    44 = 'cycling pants for women with padding'
    45 = 'cycling pants women winter'
    46 = 'cycling shorts for women'
    47 = 'cycling shorts women'
    48 = 'cycling skin suit'
    49 = 'cycling skirt'
    50 = 'cycling tight'
    51 = 'cycling tights'
    52 = 'cycling tights women'
    53 = 'cycling winter tights'
    54 = 'danskin compression capri'
    55 = 'dark brown capri leggings'
    56 = 'dark floral'
    57 = 'dark green opaque tights'
    58 = 'dark navy tights'
    59 = 'dark tights'
    60 = 'darkness rises game'
    61 = 'date due slips'
    62 = 'date like a spartan'
    63 = 'day 1 fitness weights'
    64 = 'days of the week pocket chart'
    65 = 'days of the week underwear ladies'
    66 = 'days of the week women'
    67 = 'decree leggings'
    68 = 'degree men overtime'
    69 = 'deportivas nike mujer'
    70 = 'design leggings'
    71 = 'diamondkit stretch cotton capri crop legging tights'
    72 = 'dirt bike ridding gear'
    73 = 'dkny leggings'
    74 = 'dnamic compression speed crop'
    75 = 'do it with dan'
    76 = 'dollar sign patch'
    77 = 'double couple women leggings fashion'
    78 = 'double hip brace'
    79 = 'double hip brace for women'
    80 = 'double tummy layer'
    81 = 'down skirts for women long'
    82 = 'down wash'
    83 = 'dr skin womens'
    84 = 'dragon skin sleeves'
    85 = 'dragon workout'
    86 = 'draw muscle'
    87 = 'drawstring capri leggings'
    88 = 'drawstring exercise top'
    89 = 'drawstring running pants women'
    90 = 'drawstring workout leggings'
    91 = 'dress form xl adjustable'
    92 = 'dress pants tall length for women'
    93 = 'dress to draw on'
    94 = 'dresses sun shy'
    95 = 'dri fit capri'
    96 = 'drop cloth runners'
    97 = 'drop pants for women'
    98 = 'drskin compression pants women'
    99 = 'dry cow'
    100 = 'dry fit capri leggings'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $values[$row]
}